$d = $word.ActiveDocument

# 1. Update sequencer description: "NovaSeq 6000" -> "NovaSeq X Plus (Australian Genome Research Facility)"
#    Keep the "NovaSeq" run (and its surrounding spell-check markers) untouched and only
#    replace the " 6000" suffix, so the resulting run layout matches the real edit.
$d.Content.Find.Execute(" 6000", $true, $false, $false, $false, $false,
                         $true, 1, $false, " X Plus (Australian Genome Research Facility)", 2)

# 2. Update the cached SAVEDATE field result: "16-Sep-2024" -> "4-Mar-2025"
$d.Content.Find.Execute("16-Sep-2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "4-Mar-2025", 2)
